# Inner JOIN appears to be working now
#
# Sheet1 ("Sheet1"): add two new "joined" rows (id 6 -> Dad, id 7 -> Mom)
# Sheet2: add a "Column_3" join-success boolean column, fix "MOM" -> "Mom",
#         and add a new source row (id 5 -> Mom).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet2: rename the "MOM" entry to "Mom" (row for id = 2)
# ---------------------------------------------------------------------
$ws2.Range("B3").Value = """Mom"""

# ---------------------------------------------------------------------
# Sheet2: new Column_3 header + boolean join-result column for the
# existing rows, plus a brand new source row (id 5)
# ---------------------------------------------------------------------
$ws2.Range("D1").Value = "Column_3"

$ws2.Range("D2").Value = $true
$ws2.Range("D3").Value = $true
$ws2.Range("D4").Value = $true
$ws2.Range("D5").Value = $true

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = """Mom"""
$ws2.Range("C6").Value = 8
$ws2.Range("D6").Value = $true

# ---------------------------------------------------------------------
# Sheet1: two additional joined rows
# ---------------------------------------------------------------------
$ws1.Range("A5").Value = 6
$ws1.Range("B5").Value = $true
$ws1.Range("C5").Value = """Dad"""
$ws1.Range("D5").Value = 6

$ws1.Range("A6").Value = 7
$ws1.Range("B6").Value = $false
$ws1.Range("C6").Value = """Mom"""
$ws1.Range("D6").Value = 6

# ---------------------------------------------------------------------
# Selections / active sheet: Sheet2 becomes the active tab, with its
# new Column_3 range selected; Sheet1 keeps a plain selection.
# ---------------------------------------------------------------------
$ws1.Range("D7").Select()
$ws2.Activate()
$ws2.Range("D2:D6").Select()

# Best-effort window repositioning (cosmetic; harmless if unsupported).
try {
    $win = $excel.ActiveWindow
    $win.Left = 5000
    $win.Top = 800
} catch {
}
